$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

$ws.Range("I1").Value = "specified_breakdowns"
$ws.Range("I2").Value = "T"
$ws.Range("G2").Value = "F"

$ws.Range("I1").Select()
